$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),`n                ('model',`n                 BaggingClassifier(estimator=LogisticRegression(C=3,`n                                                                max_iter=1000,`n                                                                penalty='l1',`n                                                                random_state=42,`n                                                                solver='saga'),`n                                   n_estimators=5, random_state=42))])"
$ws.Range("B2").Value = 0.6476190476190476
$ws.Range("C2").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__n_estimators': 5, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'l1', 'model__estimator__class_weight': None, 'model__estimator__C': 3}"
$ws.Range("D2").Value = 0.5000000000000001
$ws.Range("E2").Value = "[1 0 0 1 0 0 1 1 0 1 0 0]"
$ws.Range("F2").Value = "[1 1 1 0 1 1 1 1 1 1 1 1]"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.6727591036414566
$ws.Range("I2").Value = 0.02986056186189792
$ws.Range("J2").Value = 0.5721288515406163
$ws.Range("K2").Value = 0.06250057329580458

$ws.Range("A3").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',`n                                                     random_state=42))),`n                ('model',`n                 BaggingClassifier(estimator=LogisticRegression(C=0.0001,`n                                                                max_iter=1000,`n                                                                random_state=42,`n                                                                solver='saga'),`n                                   random_state=42))])"
$ws.Range("B3").Value = 0.6571428571428571
$ws.Range("C3").Value = "{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': None, 'model__n_estimators': 10, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'l2', 'model__estimator__class_weight': None, 'model__estimator__C': 0.0001}"
$ws.Range("D3").Value = 0.5333333333333333
$ws.Range("E3").Value = "[1 0 1 0 0 0 0 1 1 0 1 1]"
$ws.Range("F3").Value = "[1 1 1 1 1 0 1 1 0 1 1 0]"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.6762557077625573
$ws.Range("I3").Value = 0.03192653273908408
$ws.Range("J3").Value = 0.5598173515981734
$ws.Range("K3").Value = 0.05330226439443557

$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),`n                ('model',`n                 BaggingClassifier(estimator=LogisticRegression(C=0.001,`n                                                                max_iter=1000,`n                                                                random_state=42,`n                                                                solver='liblinear'),`n                                   random_state=42))])"
$ws.Range("B4").Value = 0.6095238095238095
$ws.Range("C4").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__n_estimators': 10, 'model__estimator__solver': 'liblinear', 'model__estimator__penalty': 'l2', 'model__estimator__class_weight': None, 'model__estimator__C': 0.001}"
$ws.Range("D4").Value = 0.823529411764706
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[1 0 1 1 1 1 1 0 1 1 0 1]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.6504960317460317
$ws.Range("I4").Value = 0.0290924160569398
$ws.Range("J4").Value = 0.5370370370370369
$ws.Range("K4").Value = 0.07299681919941778
